$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: header row (copy of row 14)
$ws.Range("C24").Value = "None"
$ws.Range("D24").Value = "D0"
$ws.Range("E24").Value = "D1"
$ws.Range("F24").Value = "D2"
$ws.Range("G24").Value = "D3"
$ws.Range("H24").Value = "D4"
$ws.Range("J24").Value = "totaldrought"

# Row 25: data row
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 100
$ws.Range("G25").Value = 100
$ws.Range("H25").Value = 74.5
$ws.Range("J25").Formula = "=100-C25"

# Row 27: header row (copy of row 17)
$ws.Range("C27").Value = "None"
$ws.Range("D27").Value = "D0"
$ws.Range("E27").Value = "D1"
$ws.Range("F27").Value = "D2"
$ws.Range("G27").Value = "D3"
$ws.Range("H27").Value = "D4"
$ws.Range("J27").Value = "sum total"

# Row 28: calculated row
$ws.Range("C28").Formula = "=C25"
$ws.Range("D28").Formula = "=D25-E25"
$ws.Range("E28").Formula = "=E25-F25"
$ws.Range("F28").Formula = "=F25-G25"
$ws.Range("G28").Formula = "=G25-H25"
$ws.Range("H28").Formula = "=H25"
$ws.Range("J28").Formula = "=SUM(C28:H28)"

# Update selection to match target state
$ws.Range("C26").Select()
